$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(4, 1).Value = 3

    $ws.Cells.Item(4, 2).NumberFormat = "@"
    $ws.Cells.Item(4, 2).Value = "2026-02-16"
    $ws.Cells.Item(4, 2).ClearFormats()

    $ws.Cells.Item(4, 3).Value = "22:57:30"
    $ws.Cells.Item(4, 4).Value = "base_strategy"
    $ws.Cells.Item(4, 5).Value = "DOWN"
    $ws.Cells.Item(4, 6).Value = 0.5
    $ws.Cells.Item(4, 7).Value = ""
    $ws.Cells.Item(4, 8).Value = "OPEN"
    $ws.Cells.Item(4, 9).Value = 0
    $ws.Cells.Item(4, 10).Value = 0
    $ws.Cells.Item(4, 11).Value = 100
    $ws.Cells.Item(4, 12).Value = 0
    $ws.Cells.Item(4, 13).Value = 0
    $ws.Cells.Item(4, 14).Value = 0.6
    $ws.Cells.Item(4, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(4, 16).Value = ""
    $ws.Cells.Item(4, 17).Value = 0
}
